# Update crypto price/volume data per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.779.71'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.20%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.306.52'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +5.88%  '

$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '602.72'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.30%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.34'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.52%  '

$ws.Range('E7').Value = '  -0.09%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.305.94'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.04%  '

$ws.Range('E10').Value = '  +2.78%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.51'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.06%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.472'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.64%  '

$ws.Range('E13').Value = '  +1.08%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.67'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.58%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.848.49'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.85%  '

$ws.Range('E16').Value = '  +0.63%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.309.06'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.97%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.865.48'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.34%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.89'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.28%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '480.31'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.70%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.16'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.55%  '

$ws.Range('E22').Value = '  +5.36%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.01'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.51%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.02'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.12%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.43'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.55%  '

$ws.Range('E26').Value = '  +0.01%  '

$ws.Range('E27').Value = '  +1.69%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.28'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.06%  '

$ws.Range('E29').Value = '  -0.15%  '

$ws.Range('E30').Value = '  +3.12%  '

$ws.Range('E31').Value = '  +4.20%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '29.37'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +10.07%  '

$ws.Range('E33').Value = '  -1.37%  '

$ws.Range('E34').Value = '  +0.78%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.09'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.86%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.97'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.97%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.88'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.83%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0748'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.25%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0401'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.00%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '427.62'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.22%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.039.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.06%  '

$ws.Range('E42').Value = '  +2.22%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.74'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.63%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.111'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.64%  '

$ws.Range('E45').Value = '  +0.20%  '

$ws.Range('E46').Value = '  +4.00%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.35'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.73%  '

$ws.Range('E48').Value = '  +0.02%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.47'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +12.29%  '

$ws.Range('E50').Value = '  +2.00%  '

$ws.Range('E51').Value = '  +2.86%  '
